$F_new = @(
    "appendText(var,appendWith)",
    "assertArrayContain(array,expected)",
    "assertArrayEqual(array1,array2,exactOrder)",
    "assertArrayNotContain(array,unexpected)",
    "assertContains(text,substring)",
    "assertCount(text,regex,expects)",
    "assertEmpty(text)",
    "assertEndsWith(text,suffix)",
    "assertEqual(expected,actual)",
    "assertNotContain(text,substring)",
    "assertNotEmpty(text)",
    "assertNotEqual(expected,actual)",
    "assertStartsWith(text,prefix)",
    "assertTextOrder(var,descending)",
    "assertVarNotPresent(var)",
    "assertVarPresent(var)",
    "clear(vars)",
    "failImmediate(text)",
    "incrementChar(var,amount,config)",
    "macro(file,sheet,name)",
    "outputToCloud(resource)",
    "prependText(var,prependWith)",
    "repeatUntil(steps,maxWaitMs)",
    "save(var,value)",
    "saveCount(text,regex,saveVar)",
    "saveMatches(text,regex,saveVar)",
    "saveReplace(text,regex,replace,saveVar)",
    "saveVariablesByPrefix(var,prefix)",
    "saveVariablesByRegex(var,regex)",
    "section(steps)",
    "split(text,delim,saveVar)",
    "startRecording()",
    "stopRecording()",
    "substringAfter(text,delim,saveVar)",
    "substringBefore(text,delim,saveVar)",
    "substringBetween(text,start,end,saveVar)",
    "verbose(text)",
    "waitFor(waitMs)"
)

$J_new = @(
    "runJUnit(className)",
    "runProgram(programPathAndParams)",
    "runProgramNoWait(programPathAndParams)",
    "tail(id,file)",
    "terminate(programName)"
)

$L_new = @(
    "assertEqual(expected,actual)",
    "assertNotEqual(expected,actual)",
    "assertPath(path)",
    "assertReadableFile(file,minByte)",
    "base64(var,file)",
    "compare(expected,actual,failFast)",
    "copyFiles(source,target)",
    "copyFilesByRegex(sourceDir,regex,target)",
    "count(var,path,pattern)",
    "deleteFiles(location,recursive)",
    "deleteFilesByRegex(sourceDir,regex)",
    "filter(source,target,matchPattern)",
    "makeDirectory(source)",
    "moveFiles(source,target)",
    "moveFilesByRegex(sourceDir,regex,target)",
    "readFile(var,file)",
    "readProperty(var,file,property)",
    "rename(target,newName)",
    "saveDiff(var,expected,actual)",
    "saveFileMeta(var,file)",
    "saveMatches(var,path,fileFilter,textFilter)",
    "searchAndReplace(file,config,saveAs)",
    "unzip(zipFile,target)",
    "validate(var,profile,inputFile)",
    "writeBase64decode(encodedSource,decodedTarget,append)",
    "writeFile(file,content,append)",
    "writeFileAsIs(file,content,append)",
    "writeProperty(file,property,value)",
    "zip(filePattern,zipFile)"
)

$Z_new = @(
    "assertAndClick(locator,label)",
    "assertAttribute(locator,attrName,value)",
    "assertAttributeContain(locator,attrName,contains)",
    "assertAttributeNotContain(locator,attrName,contains)",
    "assertAttributeNotPresent(locator,attrName)",
    "assertAttributePresent(locator,attrName)",
    "assertChecked(locator)",
    "assertContainCount(locator,text,count)",
    "assertCssNotPresent(locator,property)",
    "assertCssPresent(locator,property,value)",
    "assertElementByAttributes(nameValues)",
    "assertElementByText(locator,text)",
    "assertElementCount(locator,count)",
    "assertElementNotPresent(locator)",
    "assertElementPresent(locator)",
    "assertElementsPresent(prefix)",
    "assertFocus(locator)",
    "assertFrameCount(count)",
    "assertFramePresent(frameName)",
    "assertIECompatMode()",
    "assertIENativeMode()",
    "assertLinkByLabel(label)",
    "assertMultiSelect(locator)",
    "assertNotChecked(locator)",
    "assertNotFocus(locator)",
    "assertNotText(locator,text)",
    "assertNotVisible(locator)",
    "assertOneMatch(locator)",
    "assertScrollbarHNotPresent(locator)",
    "assertScrollbarHPresent(locator)",
    "assertScrollbarVNotPresent(locator)",
    "assertScrollbarVPresent(locator)",
    "assertSingleSelect(locator)",
    "assertTable(locator,row,column,text)",
    "assertText(locator,text)",
    "assertTextContains(locator,text)",
    "assertTextCount(locator,text,count)",
    "assertTextList(locator,list,ignoreOrder)",
    "assertTextMatches(text,minMatch,scrollTo)",
    "assertTextNotContain(locator,text)",
    "assertTextNotPresent(text)",
    "assertTextOrder(locator,descending)",
    "assertTextPresent(text)",
    "assertTitle(text)",
    "assertValue(locator,value)",
    "assertValueOrder(locator,descending)",
    "assertVisible(locator)",
    "checkAll(locator)",
    "clearLocalStorage()",
    "click(locator)",
    "clickAll(locator)",
    "clickAndWait(locator,waitMs)",
    "clickByLabel(label)",
    "clickByLabelAndWait(label,waitMs)",
    "clickOffset(locator,x,y)",
    "clickWithKeys(locator,keys)",
    "close()",
    "closeAll()",
    "deselect(locator,text)",
    "deselectMulti(locator,array)",
    "dismissInvalidCert()",
    "dismissInvalidCertPopup()",
    "doubleClick(locator)",
    "doubleClickAndWait(locator,waitMs)",
    "doubleClickByLabel(label)",
    "doubleClickByLabelAndWait(label,waitMs)",
    "dragAndDrop(fromLocator,toLocator)",
    "dragTo(fromLocator,xOffset,yOffset)",
    "editLocalStorage(key,value)",
    "executeScript(var,script)",
    "focus(locator)",
    "goBack()",
    "goBackAndWait()",
    "maximizeWindow()",
    "mouseOver(locator)",
    "open(url)",
    "openAndWait(url,waitMs)",
    "openHttpBasic(url,username,password)",
    "openIgnoreTimeout(url)",
    "refresh()",
    "refreshAndWait()",
    "resizeWindow(width,height)",
    "rightClick(locator)",
    "saveAllWindowIds(var)",
    "saveAllWindowNames(var)",
    "saveAttribute(var,locator,attrName)",
    "saveAttributeList(var,locator,attrName)",
    "saveCount(var,locator)",
    "saveDivsAsCsv(headers,rows,cells,nextPage,file)",
    "saveElement(var,locator)",
    "saveElements(var,locator)",
    "saveInfiniteDivsAsCsv(config,file)",
    "saveInfiniteTableAsCsv(config,file)",
    "saveLocalStorage(var,key)",
    "saveLocation(var)",
    "savePageAs(var,sessionIdName,url)",
    "savePageAsFile(sessionIdName,url,file)",
    "saveSelectedText(var,locator)",
    "saveSelectedValue(var,locator)",
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveValue(var,locator)",
    "saveValues(var,locator)",
    "screenshot(file,locator)",
    "scrollElement(locator,xOffset,yOffset)",
    "scrollLeft(locator,pixel)",
    "scrollPage(xOffset,yOffset)",
    "scrollRight(locator,pixel)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator)",
    "unselectAllText()",
    "updateAttribute(locator,attrName,value)",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)"
)


$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Column F ("base"): write new values (row 2 onward), then clear any leftover
# trailing cell from the previous (longer) range.
for ($i = 0; $i -lt $F_new.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $F_new[$i]
}
$ws.Cells.Item(2 + $F_new.Length, 6).ClearContents()

# Column J ("external"): write new values (row 2 onward) - grew by one row.
for ($i = 0; $i -lt $J_new.Length; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $J_new[$i]
}

# Column L ("io"): write new values (row 2 onward) - grew by one row.
for ($i = 0; $i -lt $L_new.Length; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $L_new[$i]
}

# Column Z ("web"): write new values (row 2 onward) - grew by two rows.
for ($i = 0; $i -lt $Z_new.Length; $i++) {
    $ws.Cells.Item(2 + $i, 26).Value = $Z_new[$i]
}

# Update named ranges to reflect the new extents.
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("io").RefersTo = "='#system'!`$L`$2:`$L`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
